$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original formatting of the D (Price) column, force text entry while writing
# to avoid Excel auto-converting numeric-looking strings (e.g. "9.46") into numbers,
# then restore the original style so the saved file keeps the same (unstyled) cell format.
$dRange = $ws.Range("D2:D51")
$dOrigStyle = $dRange.Style
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.128.19"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.049.18"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "247.93"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "0.663"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("D7").Value = "58.10"
$ws.Range("E7").Value = "  -3.21%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("D10").Value = "0.0775"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "15.84"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "2.348.12"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "0.846"
$ws.Range("E14").Value = "  +2.89%  "
$ws.Range("D15").Value = "5.71"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "2.049.33"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "17.98"
$ws.Range("E17").Value = "  +14.87%  "
$ws.Range("D18").Value = "37.168.75"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "74.81"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -4.01%  "
$ws.Range("D21").Value = "5.33"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").Value = "236.62"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "9.46"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  -5.64%  "
$ws.Range("D27").Value = "169.25"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "20.02"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "4.79"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "0.0616"
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("D33").Value = "4.47"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").Value = "0.0895"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "2.25"
$ws.Range("E36").Value = "  -3.04%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("E38").Value = "  +14.99%  "
$ws.Range("D39").Value = "1.33"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("D40").Value = "5.20"
$ws.Range("E40").Value = "  +14.85%  "
$ws.Range("D41").Value = "0.100"
$ws.Range("E41").Value = "  -14.36%  "
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("D43").Value = "17.14"
$ws.Range("E43").Value = "  -4.67%  "
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").Value = "95.59"
$ws.Range("E45").Value = "  -3.49%  "
$ws.Range("D46").Value = "2.44"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").Value = "1.274.51"
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("D49").Value = "6.80"
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("D50").Value = "2.233.54"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "43.70"
$ws.Range("E51").Value = "  -0.40%  "

# Restore original style on column D so number formatting / style indices match the source file
$dRange.Style = $dOrigStyle

